$wb = $excel.ActiveWorkbook

# --- db_node sheet: fix typo "soure|target" -> "source|target" ---
$dbNode = $wb.Worksheets.Item("db_node")
$dbNode.Range("E2").Value = "source|target"

# --- msq_rule sheet: update rule #1 thread counts, and rule #2 db/table mapping ---
$msqRule = $wb.Worksheets.Item("msq_rule")
$msqRule.Range("J2").Value = 2
$msqRule.Range("K2").Value = 2

$msqRule.Range("B4").Value = "msq_test"
$msqRule.Range("E4").Value = "manga"
$msqRule.Range("F4").Value = "manga"
$msqRule.Range("G4").Value = ""
$msqRule.Range("H4").Value = ""

$msqRule.Range("E5").Value = "portfolio"
$msqRule.Range("F5").Value = "portfolio"

# --- selections / active sheet, matching the final authored view state ---
$msqRule.Range("E10").Select()
$dbNode.Activate()
$dbNode.Range("C13").Select()
